$d = $word.ActiveDocument
$d.Content.Find.Execute("95-17=", $true, $false, $false, $false, $false, $true, 1, $false, "66-27=", 2) | Out-Null
$d.Content.Find.Execute("25+38=", $true, $false, $false, $false, $false, $true, 1, $false, "74-65=", 2) | Out-Null
$d.Content.Find.Execute("38+6=", $true, $false, $false, $false, $false, $true, 1, $false, "80-64=", 2) | Out-Null
$d.Content.Find.Execute("80-49=", $true, $false, $false, $false, $false, $true, 1, $false, "78+3=", 2) | Out-Null
$d.Content.Find.Execute("87-49=", $true, $false, $false, $false, $false, $true, 1, $false, "98-79=", 2) | Out-Null
$d.Content.Find.Execute("44+37=", $true, $false, $false, $false, $false, $true, 1, $false, "75-7=", 2) | Out-Null
$d.Content.Find.Execute("58+8=", $true, $false, $false, $false, $false, $true, 1, $false, "50-8=", 2) | Out-Null
$d.Content.Find.Execute("90-24=", $true, $false, $false, $false, $false, $true, 1, $false, "7+44=", 2) | Out-Null
$d.Content.Find.Execute("28+33=", $true, $false, $false, $false, $false, $true, 1, $false, "60-16=", 2) | Out-Null
$d.Content.Find.Execute("27+4=", $true, $false, $false, $false, $false, $true, 1, $false, "57+14=", 2) | Out-Null
$d.Content.Find.Execute("83-65=", $true, $false, $false, $false, $false, $true, 1, $false, "6+46=", 2) | Out-Null
$d.Content.Find.Execute("35+47=", $true, $false, $false, $false, $false, $true, 1, $false, "3+9=", 2) | Out-Null
$d.Content.Find.Execute("23-8=", $true, $false, $false, $false, $false, $true, 1, $false, "70-24=", 2) | Out-Null
$d.Content.Find.Execute("83-54=", $true, $false, $false, $false, $false, $true, 1, $false, "60-55=", 2) | Out-Null
$d.Content.Find.Execute("82-15=", $true, $false, $false, $false, $false, $true, 1, $false, "9+38=", 2) | Out-Null
$d.Content.Find.Execute("83-57=", $true, $false, $false, $false, $false, $true, 1, $false, "54-35=", 2) | Out-Null
$d.Content.Find.Execute("48+19=", $true, $false, $false, $false, $false, $true, 1, $false, "70-66=", 2) | Out-Null
$d.Content.Find.Execute("90-88=", $true, $false, $false, $false, $false, $true, 1, $false, "60-48=", 2) | Out-Null
$d.Content.Find.Execute("82-23=", $true, $false, $false, $false, $false, $true, 1, $false, "47+48=", 2) | Out-Null
$d.Content.Find.Execute("55+29=", $true, $false, $false, $false, $false, $true, 1, $false, "70-5=", 2) | Out-Null
$d.Content.Find.Execute("20-16=", $true, $false, $false, $false, $false, $true, 1, $false, "54-49=", 2) | Out-Null
$d.Content.Find.Execute("36+25=", $true, $false, $false, $false, $false, $true, 1, $false, "52+19=", 2) | Out-Null
$d.Content.Find.Execute("78+13=", $true, $false, $false, $false, $false, $true, 1, $false, "13+59=", 2) | Out-Null
$d.Content.Find.Execute("94-5=", $true, $false, $false, $false, $false, $true, 1, $false, "80-4=", 2) | Out-Null
$d.Content.Find.Execute("82-63=", $true, $false, $false, $false, $false, $true, 1, $false, "43-37=", 2) | Out-Null
$d.Content.Find.Execute("36+19=", $true, $false, $false, $false, $false, $true, 1, $false, "94-77=", 2) | Out-Null
$d.Content.Find.Execute("13+78=", $true, $false, $false, $false, $false, $true, 1, $false, "7+65=", 2) | Out-Null
$d.Content.Find.Execute("58+6=", $true, $false, $false, $false, $false, $true, 1, $false, "22-8=", 2) | Out-Null
$d.Content.Find.Execute("37+36=", $true, $false, $false, $false, $false, $true, 1, $false, "55+27=", 2) | Out-Null
$d.Content.Find.Execute("83-59=", $true, $false, $false, $false, $false, $true, 1, $false, "19+17=", 2) | Out-Null
$d.Content.Find.Execute("5+78=", $true, $false, $false, $false, $false, $true, 1, $false, "73+9=", 2) | Out-Null
$d.Content.Find.Execute("65-17=", $true, $false, $false, $false, $false, $true, 1, $false, "82-7=", 2) | Out-Null
$d.Content.Find.Execute("70-9=", $true, $false, $false, $false, $false, $true, 1, $false, "43-36=", 2) | Out-Null
$d.Content.Find.Execute("73-16=", $true, $false, $false, $false, $false, $true, 1, $false, "69+29=", 2) | Out-Null
$d.Content.Find.Execute("62-25=", $true, $false, $false, $false, $false, $true, 1, $false, "66-37=", 2) | Out-Null
$d.Content.Find.Execute("9+59=", $true, $false, $false, $false, $false, $true, 1, $false, "82-26=", 2) | Out-Null
$d.Content.Find.Execute("81-72=", $true, $false, $false, $false, $false, $true, 1, $false, "35-28=", 2) | Out-Null
$d.Content.Find.Execute("29+18=", $true, $false, $false, $false, $false, $true, 1, $false, "38+16=", 2) | Out-Null
$d.Content.Find.Execute("28+23=", $true, $false, $false, $false, $false, $true, 1, $false, "85-16=", 2) | Out-Null
$d.Content.Find.Execute("49+12=", $true, $false, $false, $false, $false, $true, 1, $false, "86+9=", 2) | Out-Null
$d.Content.Find.Execute("37+59=", $true, $false, $false, $false, $false, $true, 1, $false, "45-19=", 2) | Out-Null
$d.Content.Find.Execute("39+8=", $true, $false, $false, $false, $false, $true, 1, $false, "44-15=", 2) | Out-Null
$d.Content.Find.Execute("90-89=", $true, $false, $false, $false, $false, $true, 1, $false, "53-27=", 2) | Out-Null
$d.Content.Find.Execute("82-65=", $true, $false, $false, $false, $false, $true, 1, $false, "43+9=", 2) | Out-Null
$d.Content.Find.Execute("52-39=", $true, $false, $false, $false, $false, $true, 1, $false, "14+39=", 2) | Out-Null
$d.Content.Find.Execute("24-16=", $true, $false, $false, $false, $false, $true, 1, $false, "28+6=", 2) | Out-Null
$d.Content.Find.Execute("80-51=", $true, $false, $false, $false, $false, $true, 1, $false, "63-45=", 2) | Out-Null
$d.Content.Find.Execute("17+4=", $true, $false, $false, $false, $false, $true, 1, $false, "39+14=", 2) | Out-Null
$d.Content.Find.Execute("28+3=", $true, $false, $false, $false, $false, $true, 1, $false, "57+29=", 2) | Out-Null
$d.Content.Find.Execute("69+12=", $true, $false, $false, $false, $false, $true, 1, $false, "62-26=", 2) | Out-Null
$d.Content.Find.Execute("65+29=", $true, $false, $false, $false, $false, $true, 1, $false, "17+58=", 2) | Out-Null
$d.Content.Find.Execute("85-27=", $true, $false, $false, $false, $false, $true, 1, $false, "67+16=", 2) | Out-Null
$d.Content.Find.Execute("34-28=", $true, $false, $false, $false, $false, $true, 1, $false, "21-5=", 2) | Out-Null
$d.Content.Find.Execute("16+37=", $true, $false, $false, $false, $false, $true, 1, $false, "82-77=", 2) | Out-Null
$d.Content.Find.Execute("58+34=", $true, $false, $false, $false, $false, $true, 1, $false, "10-7=", 2) | Out-Null
$d.Content.Find.Execute("55+37=", $true, $false, $false, $false, $false, $true, 1, $false, "40-9=", 2) | Out-Null
$d.Content.Find.Execute("4+38=", $true, $false, $false, $false, $false, $true, 1, $false, "81-43=", 2) | Out-Null
$d.Content.Find.Execute("52-25=", $true, $false, $false, $false, $false, $true, 1, $false, "71-44=", 2) | Out-Null
$d.Content.Find.Execute("64+17=", $true, $false, $false, $false, $false, $true, 1, $false, "45-17=", 2) | Out-Null
$d.Content.Find.Execute("81-52=", $true, $false, $false, $false, $false, $true, 1, $false, "58+25=", 2) | Out-Null
$d.Content.Find.Execute("57+9=", $true, $false, $false, $false, $false, $true, 1, $false, "19+34=", 2) | Out-Null
$d.Content.Find.Execute("60-54=", $true, $false, $false, $false, $false, $true, 1, $false, "17+39=", 2) | Out-Null
$d.Content.Find.Execute("61-6=", $true, $false, $false, $false, $false, $true, 1, $false, "46+26=", 2) | Out-Null
$d.Content.Find.Execute("48+8=", $true, $false, $false, $false, $false, $true, 1, $false, "42-37=", 2) | Out-Null
$d.Content.Find.Execute("5+77=", $true, $false, $false, $false, $false, $true, 1, $false, "13+48=", 2) | Out-Null
$d.Content.Find.Execute("26+6=", $true, $false, $false, $false, $false, $true, 1, $false, "74-68=", 2) | Out-Null
$d.Content.Find.Execute("48-39=", $true, $false, $false, $false, $false, $true, 1, $false, "75-56=", 2) | Out-Null
$d.Content.Find.Execute("27+49=", $true, $false, $false, $false, $false, $true, 1, $false, "28+17=", 2) | Out-Null
$d.Content.Find.Execute("41-3=", $true, $false, $false, $false, $false, $true, 1, $false, "31-16=", 2) | Out-Null
$d.Content.Find.Execute("98-69=", $true, $false, $false, $false, $false, $true, 1, $false, "74+8=", 2) | Out-Null
$d.Content.Find.Execute("5+29=", $true, $false, $false, $false, $false, $true, 1, $false, "76+18=", 2) | Out-Null
$d.Content.Find.Execute("28+8=", $true, $false, $false, $false, $false, $true, 1, $false, "93-78=", 2) | Out-Null
$d.Content.Find.Execute("68-49=", $true, $false, $false, $false, $false, $true, 1, $false, "84-27=", 2) | Out-Null
$d.Content.Find.Execute("18+18=", $true, $false, $false, $false, $false, $true, 1, $false, "81-35=", 2) | Out-Null
$d.Content.Find.Execute("78-19=", $true, $false, $false, $false, $false, $true, 1, $false, "47+34=", 2) | Out-Null
$d.Content.Find.Execute("20-4=", $true, $false, $false, $false, $false, $true, 1, $false, "43-38=", 2) | Out-Null
$d.Content.Find.Execute("29+36=", $true, $false, $false, $false, $false, $true, 1, $false, "12+39=", 2) | Out-Null
$d.Content.Find.Execute("91-62=", $true, $false, $false, $false, $false, $true, 1, $false, "72-15=", 2) | Out-Null
$d.Content.Find.Execute("82-8=", $true, $false, $false, $false, $false, $true, 1, $false, "32-16=", 2) | Out-Null
$d.Content.Find.Execute("12-4=", $true, $false, $false, $false, $false, $true, 1, $false, "42-9=", 2) | Out-Null
$d.Content.Find.Execute("85-36=", $true, $false, $false, $false, $false, $true, 1, $false, "49+44=", 2) | Out-Null
$d.Content.Find.Execute("30-5=", $true, $false, $false, $false, $false, $true, 1, $false, "24-8=", 2) | Out-Null
$d.Content.Find.Execute("27+65=", $true, $false, $false, $false, $false, $true, 1, $false, "32-6=", 2) | Out-Null
$d.Content.Find.Execute("18+25=", $true, $false, $false, $false, $false, $true, 1, $false, "45+37=", 2) | Out-Null
$d.Content.Find.Execute("35-19=", $true, $false, $false, $false, $false, $true, 1, $false, "46+37=", 2) | Out-Null
$d.Content.Find.Execute("86-18=", $true, $false, $false, $false, $false, $true, 1, $false, "41-15=", 2) | Out-Null
$d.Content.Find.Execute("28+47=", $true, $false, $false, $false, $false, $true, 1, $false, "57+28=", 2) | Out-Null
$d.Content.Find.Execute("86+5=", $true, $false, $false, $false, $false, $true, 1, $false, "6+8=", 2) | Out-Null
$d.Content.Find.Execute("49+19=", $true, $false, $false, $false, $false, $true, 1, $false, "54-39=", 2) | Out-Null
$d.Content.Find.Execute("56-47=", $true, $false, $false, $false, $false, $true, 1, $false, "70-53=", 2) | Out-Null
$d.Content.Find.Execute("67-58=", $true, $false, $false, $false, $false, $true, 1, $false, "61-28=", 2) | Out-Null
$d.Content.Find.Execute("51-8=", $true, $false, $false, $false, $false, $true, 1, $false, "60-25=", 2) | Out-Null
$d.Content.Find.Execute("78+19=", $true, $false, $false, $false, $false, $true, 1, $false, "67-9=", 2) | Out-Null
$d.Content.Find.Execute("57-8=", $true, $false, $false, $false, $false, $true, 1, $false, "50-16=", 2) | Out-Null
$d.Content.Find.Execute("54+9=", $true, $false, $false, $false, $false, $true, 1, $false, "51-47=", 2) | Out-Null
$d.Content.Find.Execute("62-39=", $true, $false, $false, $false, $false, $true, 1, $false, "28+53=", 2) | Out-Null
$d.Content.Find.Execute("63-19=", $true, $false, $false, $false, $false, $true, 1, $false, "72-15=", 2) | Out-Null
$d.Content.Find.Execute("87-48=", $true, $false, $false, $false, $false, $true, 1, $false, "57-28=", 2) | Out-Null
$d.Content.Find.Execute("66-57=", $true, $false, $false, $false, $false, $true, 1, $false, "8+43=", 2) | Out-Null
$d.Content.Find.Execute("85-39=", $true, $false, $false, $false, $false, $true, 1, $false, "75-19=", 2) | Out-Null
